$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$words = @(
    "incredulous",
    "inquisitive",
    "intelligible",
    "invidious",
    "irreplaceable",
    "irrevocable",
    "jurisdiction",
    "juxtaposition",
    "kootchar",
    "languid",
    "lineage",
    "lithograph",
    "lyricist",
    "macabre",
    "madeleine",
    "maelstrom",
    "marshmallow",
    "mathematician",
    "mayonnaise",
    "memorabilia",
    "meritorious",
    "meteorology",
    "Michelangelo",
    "militia",
    "millimetre",
    "misalign",
    "misappropriate",
    "molecular",
    "nanotechnology",
    "necessitate",
    "noctarium",
    "notoriety",
    "noxious",
    "nuance",
    "occurrence",
    "oondoroo",
    "opulently",
    "orchestra",
    "ornithology",
    "parachute",
    "parliamentarian",
    "parquetry",
    "penitentiary",
    "pentameter",
    "perpendicular",
    "phalanger",
    "physician",
    "physiotherapy",
    "placebo",
    "plebiscite",
    "plesiosaur",
    "pneumatic",
    "pneumonia",
    "poinciana",
    "poinsettia",
    "portmanteau",
    "poultice",
    "precocious",
    "prejudice",
    "presbytery",
    "quarantine",
    "querulous",
    "reciprocal",
    "reconciliation",
    "repertoire",
    "secede",
    "seismology",
    "septennial",
    "seriation",
    "significance",
    "sirocco",
    "spinneret",
    "statutory",
    "stratosphere",
    "succinct",
    "suggestible",
    "supplementary",
    "surreptitious",
    "surrogate",
    "tarantella",
    "temperature",
    "therapeutic",
    "thermodynamics",
    "transcendent",
    "troposphere",
    "unconstitutional",
    "utilitarian",
    "valedictory",
    "vaudeville",
    "velocity",
    "vernacular",
    "vigilance",
    "vindictive",
    "viscount",
    "Zimbabwe",
    "aberration",
    "accelerando",
    "agglomeration",
    "algorithm",
    "amphitheatre",
    "androgynous",
    "anemometer",
    "angophora",
    "ankylosaurus",
    "anthropomorphic",
    "appendicitis",
    "arabesque",
    "arachnophobia",
    "Archimedes",
    "archipelago",
    "arpeggio",
    "assiduous",
    "Attenborough",
    "autocracy",
    "avaricious",
    "bellwether",
    "beryllium",
    "biennale",
    "bocconcini",
    "brontosaurus",
    "brucellosis",
    "bruschetta",
    "bureaucracy",
    "calicivirus",
    "capriccio",
    "caricature",
    "casuarina",
    "cataclysm",
    "caterwaul",
    "caucus",
    "celeriac",
    "centrifuge",
    "cerebellum",
    "cetacean",
    "chauffeur",
    "chicanery",
    "chlorophyll",
    "choreography",
    "chrysalis",
    "circuitous",
    "claustrophobia",
    "coccyx",
    "coefficient",
    "conjunctivitis",
    "Cotopaxi",
    "cryptozoology",
    "crystallisation",
    "cumulonimbus",
    "curmudgeonly",
    "cybersecurity",
    "decoupage",
    "deleterious",
    "dreadnought",
    "drosophila",
    "echelon",
    "effervesce",
    "electrolysis",
    "entrepreneur",
    "epitome",
    "euphoria",
    "exhilarating",
    "frugivorous",
    "fusillade",
    "hydrangea",
    "infinitesimal",
    "isthmus",
    "Kilimanjaro",
    "logarithm",
    "maestro",
    "Mendelssohn",
    "metallurgist",
    "metamorphosis",
    "Mississippi",
    "Murrumbidgee",
    "obsequious",
    "oratorio",
    "oscilloscope",
    "osteoarthritis",
    "parallelogram",
    "peninsula",
    "Philippines",
    "phosphorescence",
    "pianissimo",
    "polystyrene",
    "posthumous",
    "Pythagoras",
    "quorum",
    "ranunculus",
    "rhododendron",
    "rhythmical",
    "sclerophyll",
    "sebaceous",
    "slough",
    "styrofoam",
    "subpoena",
    "synagogue",
    "synapse",
    "tachometer",
    "thylacine",
    "triptych",
    "troubadour",
    "underprivileged",
    "vehement",
    "Venezuela",
    "vignette",
    "wurlitzer"
)

$startRow = 139
$row = $startRow
$col = 1

foreach ($word in $words) {
    $ws.Cells.Item($row, $col).Value = $word
    $col = $col + 1
    if ($col -gt 5) {
        $col = 1
        $row = $row + 1
    }
}

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 150
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K187").Select()
